$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the "Förändrad" date column (C2:C6) from 2023-09-15 (45184)
# to 2023-09-16 (45185) for each logged row.
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 3).Value = 45185
}
